# Apply "Update latest output (run 251)" changes to optimisation_result.xlsx
$wb = $excel.ActiveWorkbook

# ---- Sheet "Schedule" ----
$ws1 = $wb.Worksheets.Item("Schedule")

# Row 2: recalculated Cost ($) and Unit Cost ($/ML)
$ws1.Cells.Item(2, 5).Value = 502.96368525
$ws1.Cells.Item(2, 6).Value = 12.09628872655123

# Row 3: newly scheduled pump run
$ws1.Cells.Item(3, 1).Value = 46075.27083333334
$ws1.Cells.Item(3, 2).Value = 46075.72916666666
$ws1.Cells.Item(3, 3).Value = 11
$ws1.Cells.Item(3, 4).Value = 41.58
$ws1.Cells.Item(3, 5).Value = 945.3714269999999
$ws1.Cells.Item(3, 6).Value = 22.73620555555555
$ws1.Cells.Item(3, 1).NumberFormat = $ws1.Cells.Item(2, 1).NumberFormat
$ws1.Cells.Item(3, 2).NumberFormat = $ws1.Cells.Item(2, 2).NumberFormat

# ---- Sheet "Detailed" (sheet2): updates to existing rows 13-48 ----
$ws2 = $wb.Worksheets.Item("Detailed")

$updates = @(
  @(13, 2, 96.91313),
  @(14, 2, 98.98903),
  @(15, 2, 99.02897),
  @(16, 2, 73.2),
  @(16, 3, "historical"),
  @(17, 2, 12.48561),
  @(17, 3, "historical"),
  @(18, 2, 1.62799),
  @(18, 3, "historical"),
  @(19, 2, 0.80338),
  @(19, 3, "historical"),
  @(20, 2, 0.51),
  @(20, 3, "historical"),
  @(21, 2, 0.51),
  @(21, 3, "historical"),
  @(22, 2, 0.51),
  @(22, 3, "historical"),
  @(23, 2, -0.10684),
  @(23, 3, "historical"),
  @(24, 2, 0.01087),
  @(24, 3, "historical"),
  @(25, 3, "historical"),
  @(26, 2, 0.51),
  @(26, 3, "historical"),
  @(27, 2, 0.51),
  @(27, 3, "historical"),
  @(28, 2, 0.51),
  @(28, 3, "historical"),
  @(29, 3, "historical"),
  @(30, 2, 10.31271),
  @(30, 3, "historical"),
  @(31, 2, 57.06),
  @(31, 3, "historical"),
  @(32, 2, 57.06),
  @(32, 3, "historical"),
  @(33, 2, 50.61736),
  @(34, 2, 57.06),
  @(35, 2, 57.08),
  @(36, 2, 57.08),
  @(37, 2, 77.48911),
  @(39, 2, 106.51716),
  @(40, 2, 123.64113),
  @(41, 2, 121.01513),
  @(42, 2, 117.88908),
  @(43, 2, 136.12043),
  @(44, 2, 111.28733),
  @(45, 2, 105.79),
  @(46, 2, 103.37662),
  @(47, 2, 96.98926),
  @(48, 2, 85.45676),
)
foreach ($u in $updates) {
  $ws2.Cells.Item($u[0], $u[1]).Value = $u[2]
}

# ---- Sheet "Detailed" (sheet2): new rows 50-97 ----
$newRows = @(
  @(50, 46075, 84.79000000000001, "forecast", 46075, "OFF"),
  @(51, 46075.02083333334, 79.95, "forecast", 46075, "OFF"),
  @(52, 46075.04166666666, 79.95, "forecast", 46075, "OFF"),
  @(53, 46075.0625, 79.44889999999999, "forecast", 46075, "OFF"),
  @(54, 46075.08333333334, 79.95, "forecast", 46075, "OFF"),
  @(55, 46075.10416666666, 79.95, "forecast", 46075, "OFF"),
  @(56, 46075.125, 79.95, "forecast", 46075, "OFF"),
  @(57, 46075.14583333334, 79.95, "forecast", 46075, "OFF"),
  @(58, 46075.16666666666, 79.95, "forecast", 46075, "OFF"),
  @(59, 46075.1875, 79.95, "forecast", 46075, "OFF"),
  @(60, 46075.20833333334, 79.95, "forecast", 46075, "OFF"),
  @(61, 46075.22916666666, 79.95, "forecast", 46075, "OFF"),
  @(62, 46075.25, 79.95, "forecast", 46075, "OFF"),
  @(63, 46075.27083333334, 79.95, "forecast", 46075, "ON"),
  @(64, 46075.29166666666, 57.06, "forecast", 46075, "ON"),
  @(65, 46075.3125, 37.89, "forecast", 46075, "ON"),
  @(66, 46075.33333333334, 22.07, "forecast", 46075, "ON"),
  @(67, 46075.35416666666, 13.91287, "forecast", 46075, "ON"),
  @(68, 46075.375, 8.874409999999999, "forecast", 46075, "ON"),
  @(69, 46075.39583333334, 8.793939999999999, "forecast", 46075, "ON"),
  @(70, 46075.41666666666, 12.4639, "forecast", 46075, "ON"),
  @(71, 46075.4375, 22.07, "forecast", 46075, "ON"),
  @(72, 46075.45833333334, 35.88, "forecast", 46075, "ON"),
  @(73, 46075.47916666666, 35.88, "forecast", 46075, "ON"),
  @(74, 46075.5, 35.88, "forecast", 46075, "ON"),
  @(75, 46075.52083333334, 37.89, "forecast", 46075, "ON"),
  @(76, 46075.54166666666, 37.89, "forecast", 46075, "ON"),
  @(77, 46075.5625, 64.89, "forecast", 46075, "ON"),
  @(78, 46075.58333333334, 57.06, "forecast", 46075, "ON"),
  @(79, 46075.60416666666, 64.89, "forecast", 46075, "ON"),
  @(80, 46075.625, 70.39059, "forecast", 46075, "ON"),
  @(81, 46075.64583333334, 64.38733000000001, "forecast", 46075, "ON"),
  @(82, 46075.66666666666, 64.89, "forecast", 46075, "ON"),
  @(83, 46075.6875, 64.89, "forecast", 46075, "ON"),
  @(84, 46075.70833333334, 71.70868, "forecast", 46075, "ON"),
  @(85, 46075.72916666666, 84.04841999999999, "forecast", 46075, "OFF"),
  @(86, 46075.75, 78.54272, "forecast", 46075, "OFF"),
  @(87, 46075.77083333334, 100.01, "forecast", 46075, "OFF"),
  @(88, 46075.79166666666, 98.07807, "forecast", 46075, "OFF"),
  @(89, 46075.8125, 89.01005000000001, "forecast", 46075, "OFF"),
  @(90, 46075.83333333334, 84.79000000000001, "forecast", 46075, "OFF"),
  @(91, 46075.85416666666, 84.79000000000001, "forecast", 46075, "OFF"),
  @(92, 46075.875, 78.21836, "forecast", 46075, "OFF"),
  @(93, 46075.89583333334, 75.68000000000001, "forecast", 46075, "OFF"),
  @(94, 46075.91666666666, 71.15564000000001, "forecast", 46075, "OFF"),
  @(95, 46075.9375, 68.37699000000001, "forecast", 46075, "OFF"),
  @(96, 46075.95833333334, 68.22456, "forecast", 46075, "OFF"),
  @(97, 46075.97916666666, 73.2, "forecast", 46075, "OFF"),
)
foreach ($r in $newRows) {
  $rowNum = $r[0]
  $ws2.Cells.Item($rowNum, 1).Value = $r[1]
  $ws2.Cells.Item($rowNum, 2).Value = $r[2]
  $ws2.Cells.Item($rowNum, 3).Value = $r[3]
  $ws2.Cells.Item($rowNum, 4).Value = $r[4]
  $ws2.Cells.Item($rowNum, 5).Value = $r[5]
  $ws2.Cells.Item($rowNum, 1).NumberFormat = $ws2.Cells.Item(2, 1).NumberFormat
  $ws2.Cells.Item($rowNum, 4).NumberFormat = $ws2.Cells.Item(2, 4).NumberFormat
}

Write-Output "Done: Schedule dims -> A1:F3, Detailed dims -> A1:E97"